# Applies the "map.xlsx" edit:
#  - C4:F4 keep the "SHELF_N" text, but it is renamed to "SHELF_N_CC:aa"
#  - G4:J4 become a new string "SHELF_N_CC:ab" (J4 was GROUND before)
#  - C5:F5 keep the "SHELF_S" text, but it is renamed to "SHELF_S_CC:ac"
#  - G5:J5 become a new string "SHELF_S_CC:ad" (J5 was GROUND before)
#  - a new column L is populated with "GROUND" for rows 2-27
#  - the used range grows from A1:K27 to A1:L27
#  - the active selection moves to M8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shelf relabeling on row 4 (north shelf) ---
$ws.Range("C4:F4").Value = "SHELF_N_CC:aa"
$ws.Range("G4:J4").Value = "SHELF_N_CC:ab"

# --- Shelf relabeling on row 5 (south shelf) ---
$ws.Range("C5:F5").Value = "SHELF_S_CC:ac"
$ws.Range("G5:J5").Value = "SHELF_S_CC:ad"

# --- New column L: GROUND for every row except the top wall row ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 12).Value = "GROUND"
}

# --- Move the active selection to M8, matching the saved sheetView ---
$ws.Range("M8").Select()
